$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '37.937.98'
Set-TextValue "E2" '  -0.46%  '
Set-TextValue "D3" '2.037.24'
Set-TextValue "E3" '  -0.87%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '228.08'
Set-TextValue "E5" '  -0.63%  '
Set-TextValue "E6" '  -0.79%  '
Set-TextValue "D7" '60.83'
Set-TextValue "E7" '  +3.14%  '
Set-TextValue "E8" '  +0.09%  '
Set-TextValue "E9" '  -1.82%  '
Set-TextValue "E10" '  +0.64%  '
Set-TextValue "E11" '  +0.44%  '
Set-TextValue "D12" '2.339.91'
Set-TextValue "E12" '  -0.68%  '
Set-TextValue "D13" '14.52'
Set-TextValue "E13" '  -0.90%  '
Set-TextValue "D14" '21.45'
Set-TextValue "E14" '  +2.24%  '
Set-TextValue "E15" '  +1.39%  '
Set-TextValue "D16" '5.17'
Set-TextValue "E16" '  -2.05%  '
Set-TextValue "D17" '2.046.27'
Set-TextValue "E17" '  -0.70%  '
Set-TextValue "D18" '37.892.64'
Set-TextValue "E18" '  -0.24%  '
Set-TextValue "D19" '69.82'
Set-TextValue "E19" '  +0.03%  '
Set-TextValue "D20" '5.90'
Set-TextValue "E20" '  -6.56%  '
Set-TextValue "E21" '  -1.50%  '
Set-TextValue "D22" '224.22'
Set-TextValue "E22" '  -0.14%  '
Set-TextValue "E23" '  +0.03%  '
Set-TextValue "D24" '2.42'
Set-TextValue "E24" '  -0.14%  '
Set-TextValue "E25" '  +0.01%  '
Set-TextValue "D26" '9.36'
Set-TextValue "E26" '  +0.84%  '
Set-TextValue "D27" '167.35'
Set-TextValue "E27" '  +0.66%  '
Set-TextValue "E28" '  -2.09%  '
Set-TextValue "D29" '18.89'
Set-TextValue "E29" '  -0.77%  '
Set-TextValue "E30" '  -3.53%  '
Set-TextValue "E31" '  +0.64%  '
Set-TextValue "E32" '  +9.57%  '
Set-TextValue "D33" '4.41'
Set-TextValue "E33" '  -2.94%  '
Set-TextValue "E34" '  +0.20%  '
Set-TextValue "E35" '  -1.65%  '
Set-TextValue "D36" '6.38'
Set-TextValue "E36" '  +4.93%  '
Set-TextValue "E37" '  -0.88%  '
Set-TextValue "D38" '3.37'
Set-TextValue "E38" '  +3.06%  '
Set-TextValue "E39" '  +0.03%  '
Set-TextValue "B40" 'InjectiveProtocol'
Set-TextValue "C40" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D40" '17.71'
Set-TextValue "E40" '  +5.32%  '
Set-TextValue "B41" 'Maker'
Set-TextValue "C41" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D41" '1.538.97'
Set-TextValue "E41" '  +0.24%  '
Set-TextValue "D42" '0.0218'
Set-TextValue "E42" '  +0.66%  '
Set-TextValue "D43" '96.39'
Set-TextValue "E43" '  -1.82%  '
Set-TextValue "E44" '  -2.54%  '
Set-TextValue "D45" '0.0914'
Set-TextValue "E45" '  -0.81%  '
Set-TextValue "E46" '  -2.37%  '
Set-TextValue "E47" '  -1.43%  '
Set-TextValue "E48" '  -0.83%  '
Set-TextValue "D49" '2.97'
Set-TextValue "E49" '  -0.03%  '
Set-TextValue "D50" '7.11'
Set-TextValue "E50" '  -0.31%  '
Set-TextValue "D51" '2.228.82'
